$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells stay text-formatted so values like
# "314.44" or "35.00" are not reinterpreted as numbers and keep their
# exact original formatting (e.g. trailing zeros, double-dot thousands).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.543.42"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.459.62"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.44"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.40"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +2.18%  "
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").Value = "  +4.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.42"
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0799"
$ws.Range("E11").Value = "  +3.68%  "
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.837.35"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.84"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.83"
$ws.Range("E15").Value = "  +4.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.442.98"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.769"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.526.41"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.46"
$ws.Range("E19").Value = "  +4.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0936"
$ws.Range("E20").Value = "  +3.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.79"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.34"
$ws.Range("E22").Value = "  +4.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.64"
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.70"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.25"
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.65"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.00"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.48"
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0758"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.39"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.88"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.102"
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.93"
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.969.02"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0281"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.38"
$ws.Range("E45").Value = "  -8.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.93"
$ws.Range("E47").Value = "  +3.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.695.68"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.36"
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.29"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("E51").Value = "  -1.19%  "
